$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the REQ values in A2:A6 with new data
$ws.Range("A2").Value = "REQ0717399"
$ws.Range("A3").Value = "REQ0715319"
$ws.Range("A4").Value = "REQ0714742"
$ws.Range("A5").Value = "REQ0714563"
$ws.Range("A6").Value = "REQ0714396"

# Remove the last row (row 7), which previously held REQ0561303
$ws.Range("A7").EntireRow.Delete()
